$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Insert Dish Category -> CategoryRecipes column (O) gets "C"
$ws.Range("A10").Value = "Insert Dish Category"
$ws.Range("A10").WrapText = $true
$ws.Range("O10").Value = "C"

# Row 11: Update Dish Category -> CategoryRecipes column (O) gets "U"
$ws.Range("A11").Value = "Update Dish Category"
$ws.Range("A11").WrapText = $true
$ws.Range("O11").Value = "U"

# Row 12: Delete Dish Category -> CategoryRecipes column (O) gets "D"
$ws.Range("A12").Value = "Delete Dish Category"
$ws.Range("A12").WrapText = $true
$ws.Range("O12").Value = "D"

# Row 13: Insert Cuisine Type -> Type of cuisine column (N) gets "C"
$ws.Range("A13").Value = "Insert Cuisine Type"
$ws.Range("A13").WrapText = $true
$ws.Range("N13").Value = "C"

# Row 14: Update Cuisine Type -> Type of cuisine column (N) gets "U"
$ws.Range("A14").Value = "Update Cuisine Type"
$ws.Range("A14").WrapText = $true
$ws.Range("N14").Value = "U"

# Row 15 (new row): Delete Cuisine Type -> Type of cuisine column (N) gets "D"
$ws.Range("A15").Value = "Delete Cuisine Type"
$ws.Range("A15").WrapText = $true
$ws.Range("N15").Value = "D"
$ws.Range("N9").Copy()
$ws.Range("N15").PasteSpecial(-4122)

# Update selection to match new active cell
$ws.Range("M14").Select()
